$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.383390763675381
$ws.Range("C2").Value = 6.398501588092466
$ws.Range("D2").Value = 4.516332085582391
$ws.Range("F2").Value = 20.13466825117663
$ws.Range("G2").Value = 21.52392639752492
$ws.Range("H2").Value = 12.61216335643556
$ws.Range("K2").Value = 8.377039131893014
$ws.Range("N2").Value = 17.08841876230734
$ws.Range("O2").Value = 18.20840019444088
$ws.Range("B3").Value = 8.051603992731851
$ws.Range("C3").Value = 6.313838272946482
$ws.Range("D3").Value = 4.40406399120435
$ws.Range("F3").Value = 20.19153138429924
$ws.Range("G3").Value = 21.62587227222362
$ws.Range("H3").Value = 12.65686267874412
$ws.Range("K3").Value = 8.10791070653498
$ws.Range("N3").Value = 17.12981588840008
$ws.Range("O3").Value = 18.28749447322488
$ws.Range("B4").Value = 7.841484863862886
$ws.Range("C4").Value = 6.261117228711338
$ws.Range("D4").Value = 4.332837377937192
$ws.Range("F4").Value = 20.23182984366015
$ws.Range("G4").Value = 21.6958978076013
$ws.Range("H4").Value = 12.68609766366574
$ws.Range("K4").Value = 7.936556261092459
$ws.Range("N4").Value = 17.15676634705311
$ws.Range("O4").Value = 18.33973899806802
$ws.Range("B5").Value = 7.754375333004056
$ws.Range("C5").Value = 6.239464286947438
$ws.Range("D5").Value = 4.303259176983968
$ws.Range("F5").Value = 20.24960243589841
$ws.Range("G5").Value = 21.72629322583527
$ws.Range("H5").Value = 12.69846158934152
$ws.Range("K5").Value = 7.865258215294769
$ws.Range("N5").Value = 17.16813498775783
$ws.Range("O5").Value = 18.36195385560065
$ws.Range("B6").Value = 7.739825161412315
$ws.Range("C6").Value = 6.235859139096298
$ws.Range("D6").Value = 4.298315086586202
$ws.Range("F6").Value = 20.25263503483898
$ws.Range("G6").Value = 21.73145239964096
$ws.Range("H6").Value = 12.70054182869462
$ws.Range("K6").Value = 7.853332530262238
$ws.Range("N6").Value = 17.17004608453155
$ws.Range("O6").Value = 18.36569845189284
$ws.Range("B7").Value = 7.84031590787137
$ws.Range("C7").Value = 6.260825870244934
$ws.Range("D7").Value = 4.332440681248701
$ws.Range("F7").Value = 20.23206406690626
$ws.Range("G7").Value = 21.69630021409005
$ws.Range("H7").Value = 12.68626258347852
$ws.Range("K7").Value = 7.935600569059294
$ws.Range("N7").Value = 17.15691810415304
$ws.Range("O7").Value = 18.34003485181801
$ws.Range("B8").Value = 8.270389833892219
$ws.Range("C8").Value = 6.369472676659331
$ws.Range("D8").Value = 4.478111907472025
$ws.Range("F8").Value = 20.15315528129901
$ws.Range("G8").Value = 21.5575298549145
$ws.Range("H8").Value = 12.62720456962764
$ws.Range("K8").Value = 8.285549215466204
$ws.Range("N8").Value = 17.10237488763453
$ws.Range("O8").Value = 18.23490765731063
$ws.Range("B9").Value = 9.058189613682215
$ws.Range("C9").Value = 6.575960204360687
$ws.Range("D9").Value = 4.744610859569119
$ws.Range("F9").Value = 20.04127024203244
$ws.Range("G9").Value = 21.34474581413324
$ws.Range("H9").Value = 12.52556899504772
$ws.Range("K9").Value = 8.92074673433083
$ws.Range("N9").Value = 17.00754071345407
$ws.Range("O9").Value = 18.05798441773572
$ws.Range("B10").Value = 9.597619209134322
$ws.Range("C10").Value = 6.722731277801544
$ws.Range("D10").Value = 4.927558276605956
$ws.Range("F10").Value = 19.98535551648219
$ws.Range("G10").Value = 21.22509954320446
$ws.Range("H10").Value = 12.45951046437427
$ws.Range("K10").Value = 9.353282376604049
$ws.Range("N10").Value = 16.94520994514743
$ws.Range("O10").Value = 17.94585857839899
$ws.Range("B11").Value = 9.833549983907938
$ws.Range("C11").Value = 6.788241650761391
$ws.Range("D11").Value = 5.007782134841358
$ws.Range("F11").Value = 19.96565340070045
$ws.Range("G11").Value = 21.17873641844142
$ws.Range("H11").Value = 12.43132291283813
$ws.Range("K11").Value = 9.54211538058591
$ws.Range("N11").Value = 16.91843874447682
$ws.Range("O11").Value = 17.89873655645885
$ws.Range("B12").Value = 9.921467724077813
$ws.Range("C12").Value = 6.812852996917618
$ws.Range("D12").Value = 5.037713690017422
$ws.Range("F12").Value = 19.95901875892003
$ws.Range("G12").Value = 21.16234723822602
$ws.Range("H12").Value = 12.42091647825356
$ws.Range("K12").Value = 9.612444209581177
$ws.Range("N12").Value = 16.90852814590681
$ws.Range("O12").Value = 17.88145211214852
$ws.Range("B13").Value = 9.902597326923519
$ws.Range("C13").Value = 6.807561458293019
$ws.Range("D13").Value = 5.031287536174974
$ws.Range("F13").Value = 19.96041088471695
$ws.Range("G13").Value = 21.16582489539234
$ws.Range("H13").Value = 12.42314579394798
$ws.Range("K13").Value = 9.59735058923585
$ws.Range("N13").Value = 16.91065248108211
$ws.Range("O13").Value = 17.88514971949403
$ws.Range("B14").Value = 9.840811905212227
$ws.Range("C14").Value = 6.79027045497044
$ws.Range("D14").Value = 5.010253676912436
$ws.Range("F14").Value = 19.96509100032795
$ws.Range("G14").Value = 21.17736462844183
$ws.Range("H14").Value = 12.43046140812513
$ws.Range("K14").Value = 9.547925187543315
$ws.Range("N14").Value = 16.91761884613209
$ws.Range("O14").Value = 17.89730333273044
$ws.Range("B15").Value = 9.802779343471695
$ws.Range("C15").Value = 6.779653245411568
$ws.Range("D15").Value = 4.997311103718234
$ws.Range("F15").Value = 19.96806533215681
$ws.Range("G15").Value = 21.18458531969994
$ws.Range("H15").Value = 12.43497727140415
$ws.Range("K15").Value = 9.517496230247257
$ws.Range("N15").Value = 16.92191549964104
$ws.Range("O15").Value = 17.90482068493962
$ws.Range("B16").Value = 9.582004566909129
$ws.Range("C16").Value = 6.718423412723627
$ws.Range("D16").Value = 4.922253682763991
$ws.Range("F16").Value = 19.98675862197219
$ws.Range("G16").Value = 21.22829244447822
$ws.Range("H16").Value = 12.46139004650708
$ws.Range("K16").Value = 9.340778787391731
$ws.Range("N16").Value = 16.94699131491901
$ws.Range("O16").Value = 17.94901637014111
$ws.Range("B17").Value = 9.444095341098684
$ws.Range("C17").Value = 6.68052797488379
$ws.Range("D17").Value = 4.875428109155937
$ws.Range("F17").Value = 19.99969619079611
$ws.Range("G17").Value = 21.2571766178868
$ws.Range("H17").Value = 12.47807032746573
$ws.Range("K17").Value = 9.230308946146756
$ws.Range("N17").Value = 16.96277961560602
$ws.Range("O17").Value = 17.9771246662989
$ws.Range("B18").Value = 9.363886528527932
$ws.Range("C18").Value = 6.658614130727276
$ws.Range("D18").Value = 4.848213781865559
$ws.Range("F18").Value = 20.00767715870404
$ws.Range("G18").Value = 21.27454866500735
$ws.Range("H18").Value = 12.4878397364696
$ws.Range("K18").Value = 9.166025190526764
$ws.Range("N18").Value = 16.97200969949479
$ws.Range("O18").Value = 17.9936573871604
$ws.Range("B19").Value = 9.336578949875864
$ws.Range("C19").Value = 6.651174794415401
$ws.Range("D19").Value = 4.838951646883035
$ws.Range("F19").Value = 20.01047200144282
$ws.Range("G19").Value = 21.28056059072995
$ws.Range("H19").Value = 12.49117761650933
$ws.Range("K19").Value = 9.144133190839682
$ws.Range("N19").Value = 16.97516046930259
$ws.Range("O19").Value = 17.99931783814202
$ws.Range("B20").Value = 9.458868315221606
$ws.Range("C20").Value = 6.684574271334891
$ws.Range("D20").Value = 4.880442026651627
$ws.Range("F20").Value = 19.99826310276375
$ws.Range("G20").Value = 21.25402327407497
$ws.Range("H20").Value = 12.47627653688886
$ws.Range("K20").Value = 9.242145983602033
$ws.Range("N20").Value = 16.96108350025349
$ws.Range("O20").Value = 17.97409464729494
$ws.Range("B21").Value = 9.858998896656797
$ws.Range("C21").Value = 6.795354680556211
$ws.Range("D21").Value = 5.016444096539415
$ws.Range("F21").Value = 19.96369390552621
$ws.Range("G21").Value = 21.17394338253485
$ws.Range("H21").Value = 12.42830537512481
$ws.Range("K21").Value = 9.56247487896866
$ws.Range("N21").Value = 16.91556649794867
$ws.Range("O21").Value = 17.89371832499891
$ws.Range("B22").Value = 10.11218299486678
$ws.Range("C22").Value = 6.8666069107995
$ws.Range("D22").Value = 5.102714643348321
$ws.Range("F22").Value = 19.9459166449067
$ws.Range("G22").Value = 21.12841508962916
$ws.Range("H22").Value = 12.39851298013243
$ws.Range("K22").Value = 9.76494690983831
$ws.Range("N22").Value = 16.88714179209112
$ws.Range("O22").Value = 17.84445026971054
$ws.Range("B23").Value = 9.97783439999869
$ws.Range("C23").Value = 6.828688307747862
$ws.Range("D23").Value = 5.056914604302035
$ws.Range("F23").Value = 19.9549636463572
$ws.Range("G23").Value = 21.15208894635397
$ws.Range("H23").Value = 12.41427113712148
$ws.Range("K23").Value = 9.657524753032044
$ws.Range("N23").Value = 16.90219171108004
$ws.Range("O23").Value = 17.87044668978772
$ws.Range("B24").Value = 9.452192323483963
$ws.Range("C24").Value = 6.682745335910099
$ws.Range("D24").Value = 4.878176148070371
$ws.Range("F24").Value = 19.99890931065164
$ws.Range("G24").Value = 21.25544651494096
$ws.Range("H24").Value = 12.47708694985874
$ws.Range("K24").Value = 9.236796864324884
$ws.Range("N24").Value = 16.96184983681749
$ws.Range("O24").Value = 17.97546335745997
$ws.Range("B25").Value = 8.851632941959551
$ws.Range("C25").Value = 6.520900943369861
$ws.Range("D25").Value = 4.6746943738831
$ws.Range("F25").Value = 20.06693202095495
$ws.Range("G25").Value = 21.39590238471622
$ws.Range("H25").Value = 12.55154946875415
$ws.Range("K25").Value = 8.754709077504936
$ws.Range("N25").Value = 17.03190280430455
$ws.Range("O25").Value = 18.10271342893637
